# LDLC smartphone price-tracking sheet: append a new timestamped price
# snapshot column, inserted right before the existing "nom" / "url_produit"
# columns (which shift one column to the right, GS->GT and GT->GU).
#
# Before: ... GR=<last price col>  GS=nom            GT=url_produit
# After:  ... GR=<last price col>  GS=<new timestamp col, copy of GR>  GT=nom  GU=url_produit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at GS, pushing the old GS ("nom") to GT and the old
# GT ("url_produit") to GU. Formatting of the surrounding cells is carried
# along automatically by Excel's column insert semantics.
$ws.Columns("GS").Insert()

# New column header: the timestamp of this price check.
$ws.Range("GS1").Value = "2026-02-06 08:28:21"

# The new snapshot column simply repeats the most recent price (the value
# that was already sitting in column GR) for every product row that has a
# price there. Rows below the last real product (81-210) have nothing in
# GR, so they stay blank, matching the rest of that row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $price = $ws.Cells.Item($r, 200).Value2   # column GR = 200
    if ($null -ne $price -and $price -ne "") {
        $ws.Cells.Item($r, 201).Value = $price   # column GS = 201
    }
}
